# Daily scrape update - 2025-09-13 02:58:00 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# (width stored in xlsx = ColumnWidth + 5/6; subtract that offset so the
#  saved file ends up with the exact target widths)
$ws.Columns.Item(3).ColumnWidth = 54.166666666666664   # C: 67 -> 55
$ws.Columns.Item(4).ColumnWidth = 69.16666666666667    # D: 65 -> 70
$ws.Columns.Item(7).ColumnWidth = 14.166666666666666   # G: 16 -> 15
$ws.Columns.Item(8).ColumnWidth = 22.166666666666668   # H: 50 -> 23

# --- Keep the OPPORTUNITY ID column as text (the ids are numeric-looking
#     strings, not numbers) ---
$ws.Range("A2:A13").NumberFormat = "@"

# --- Row 2 ---
$ws.Range("A2").Value = "1327553"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1327553"
$ws.Range("C2").Value = "SM illustrator"
$ws.Range("D2").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "Perfect design firm"

# --- Row 3 ---
$ws.Range("A3").Value = "1327551"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1327551"
$ws.Range("C3").Value = "Video Editor"
$ws.Range("D3").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F3").Value = "0 applicants"
$ws.Range("H3").Value = "Perfect design firm"

# --- Row 4 ---
$ws.Range("A4").Value = "1327547"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1327547"
$ws.Range("D4").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("H4").Value = "Perfect design firm"

# --- Row 5 ---
$ws.Range("A5").Value = "1327499"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327499"
$ws.Range("C5").Value = "Full Stack Developer"
$ws.Range("D5").Value = "Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"
$ws.Range("F5").Value = "0 applicants"
$ws.Range("H5").Value = "Skyline Egypt Tours"

# --- Row 6 ---
$ws.Range("A6").Value = "1327497"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1327497"
$ws.Range("C6").Value = "SEO Specialist"
$ws.Range("D6").Value = "Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"
$ws.Range("F6").Value = "0 applicants"
$ws.Range("H6").Value = "Skyline Egypt Tours"

# --- Row 7 ---
$ws.Range("A7").Value = "1327495"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1327495"
$ws.Range("C7").Value = "Content Creator"
$ws.Range("D7").Value = "Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"
$ws.Range("F7").Value = "2 applicants"
$ws.Range("G7").Value = "3 - 6 Months"
$ws.Range("H7").Value = "Skyline Egypt Tours"

# --- Row 8 ---
$ws.Range("A8").Value = "1327475"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1327475"
$ws.Range("C8").Value = "Property Consultant"
$ws.Range("D8").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Bold Routes"

# --- Row 9 ---
$ws.Range("A9").Value = "1327433"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1327433"
$ws.Range("C9").Value = "Content Creator"
$ws.Range("D9").Value = "Zagazig, El-Hariry, Zagazig 1, Al-Sharqia Governorate, Egypt"
$ws.Range("G9").Value = "9 - 12 Weeks"
$ws.Range("H9").Value = "Admixy"

# --- Row 10 ---
$ws.Range("A10").Value = "1327335"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1327335"
$ws.Range("C10").Value = "2D Animator"
$ws.Range("D10").Value = "Al Manteqah Ath Thamenah, Nasr City, Cairo Governorate, Egypt"
$ws.Range("G10").Value = "9 - 12 Weeks"
$ws.Range("H10").Value = "Sparks Studio"

# --- Row 11 ---
$ws.Range("A11").Value = "1326697"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1326697"
$ws.Range("C11").Value = "Social Media Strategist for an innovative AI product"
$ws.Range("D11").Value = "Cairo, Cairo Governorate, Egypt"
$ws.Range("F11").Value = "1 applicant"
$ws.Range("H11").Value = "Metrics"

# --- Row 12 ---
$ws.Range("A12").Value = "1326535"
$ws.Range("B12").Value = "https://aiesec.org/opportunity/global-talent/1326535"
$ws.Range("C12").Value = "ACCOUNTANT"
$ws.Range("D12").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("F12").Value = "21 applicants"
$ws.Range("G12").Value = "3 - 6 Months"
$ws.Range("H12").Value = "Egypt holiday travel"

# --- Row 13 (new row) ---
$ws.Range("A13").Value = "1326381"
$ws.Range("B13").Value = "https://aiesec.org/opportunity/global-talent/1326381"
$ws.Range("C13").Value = "Business Developer"
$ws.Range("D13").Value = "Sheraton Al Matar, El Nozha, Cairo Governorate, Egypt"
$ws.Range("E13").Value = "No"
$ws.Range("F13").Value = "1 applicant"
$ws.Range("G13").Value = "3 - 6 Months"
$ws.Range("H13").Value = "1 applicant"
